$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-30 18:28:28"

# --- Step 1: shift the existing data rows (2..5) down to (3..6), bottom-up
# so we never clobber a row before reading it. Column A (timestamp) is
# re-stamped with the new fetch time on every row; B..E, G, H are copied
# as-is. Column F (URL) and the hyperlinks are rebuilt separately below.
for ($r = 5; $r -ge 2; $r--) {
    $dst = $r + 1
    $ws.Range("B$dst").Value = $ws.Range("B$r").Value2
    $ws.Range("C$dst").Value = $ws.Range("C$r").Value2
    $ws.Range("D$dst").Value = $ws.Range("D$r").Value2
    $ws.Range("E$dst").Value = $ws.Range("E$r").Value2
    $ws.Range("F$dst").Value = $ws.Range("F$r").Value2
    $ws.Range("G$dst").Value = $ws.Range("G$r").Value2
    $ws.Range("H$dst").Value = $ws.Range("H$r").Value2
    $ws.Range("A$dst").Value = $newTimestamp
}

# --- Step 2: write the brand-new row 2 (latest scraped listing)
$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = "大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5455098"
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# --- Step 3: rebuild the hyperlinks on F2:F6 from scratch so the
# relationship ids come out in clean top-to-bottom order (rId1..rId5),
# matching each row's (now shifted) URL.
for ($r = 2; $r -le 6; $r++) {
    $target = $ws.Range("F$r")
    if ($target.Hyperlinks.Count -gt 0) {
        $target.Hyperlinks.Delete()
    }
}

$urls = @(
    "https://www.lancers.jp/work/detail/5455098",
    "https://www.lancers.jp/work/detail/5445159",
    "https://www.lancers.jp/work/detail/5445154",
    "https://www.lancers.jp/work/detail/5463183",
    "https://www.lancers.jp/work/detail/5463296"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $r = $i + 2
    $cell = $ws.Range("F$r")
    $ws.Hyperlinks.Add($cell, $urls[$i]) | Out-Null
    $cell.Style = "Hyperlink"
}
